# Implements "second scenario and 3 small test scenarios":
#  - adds a new "BRAKE" column (K) to Sheet1, header + 0 for every data row
#  - hard-codes a couple of cells in column F (rows 4, 6, 7) that used to be
#    driven by the shared decelerate formula, replacing them with plain
#    values
#  - moves the sheet's active selection/scroll position back to F7

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New "BRAKE" column (K) -------------------------------------------
$ws.Range("K1").Value = "BRAKE"
$ws.Range("K2:K156").Value = 0

# --- Hard-code the three cells that stop following the shared formula --
$ws.Range("F4").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 0

# --- Restore the scroll/selection state shown in the saved file --------
$ws.Range("F7").Select() | Out-Null
